$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.024.23'
$ws.Range('E2').Value = '  +0.43%  '
$ws.Range('D3').Value = '1.926.36'
$ws.Range('E3').Value = '  +1.15%  '
$cell = $ws.Range('D4')
$cell.Value = "'1.005"
$cell.Style = 'Normal'
$ws.Range('E4').Value = '  +0.01%  '
$cell = $ws.Range('D5')
$cell.Value = "'325.54"
$cell.Style = 'Normal'
$ws.Range('E5').Value = '  +0.44%  '
$ws.Range('E6').Value = '  +0.11%  '
$cell = $ws.Range('D7')
$cell.Value = "'0.4585"
$cell.Style = 'Normal'
$ws.Range('E7').Value = '  -0.04%  '
$cell = $ws.Range('D8')
$cell.Value = "'0.3821"
$cell.Style = 'Normal'
$ws.Range('E8').Value = '  +0.30%  '
$cell = $ws.Range('D9')
$cell.Value = "'0.07759"
$cell.Style = 'Normal'
$ws.Range('E9').Value = '  +0.38%  '
$cell = $ws.Range('D10')
$cell.Value = "'0.9798"
$cell.Style = 'Normal'
$ws.Range('E10').Value = '  -0.23%  '
$ws.Range('E11').Value = '  +2.51%  '
$ws.Range('D12').Value = '1.942.24'
$ws.Range('E12').Value = '  -0.86%  '
$cell = $ws.Range('D13')
$cell.Value = "'5.716"
$cell.Style = 'Normal'
$ws.Range('E13').Value = '  +0.71%  '
$cell = $ws.Range('D14')
$cell.Value = "'6.989"
$cell.Style = 'Normal'
$ws.Range('E14').Value = '  +0.12%  '
$ws.Range('E15').Value = '  -0.67%  '
$cell = $ws.Range('D16')
$cell.Value = "'84.89"
$cell.Style = 'Normal'
$ws.Range('E16').Value = '  +0.90%  '
$cell = $ws.Range('D17')
$cell.Value = "'1.006"
$cell.Style = 'Normal'
$ws.Range('E17').Value = '  +0.13%  '
$cell = $ws.Range('D18')
$cell.Value = "'0.000009501"
$cell.Style = 'Normal'
$ws.Range('E18').Value = '  -0.50%  '
$cell = $ws.Range('D19')
$cell.Value = "'16.74"
$cell.Style = 'Normal'
$ws.Range('E19').Value = '  -0.08%  '
$ws.Range('E20').Value = '  +0.09%  '
$ws.Range('D21').Value = '29.033.12'
$ws.Range('E21').Value = '  +0.45%  '
$cell = $ws.Range('D22')
$cell.Value = "'5.356"
$cell.Style = 'Normal'
$ws.Range('E22').Value = '  +0.44%  '
$cell = $ws.Range('D23')
$cell.Value = "'11.09"
$cell.Style = 'Normal'
$ws.Range('E23').Value = '  +1.56%  '
$ws.Range('D24').Value = '2.198.68'
$ws.Range('E24').Value = '  +1.25%  '
$cell = $ws.Range('D25')
$cell.Value = "'2.057"
$cell.Style = 'Normal'
$ws.Range('E25').Value = '  -1.18%  '
$cell = $ws.Range('D26')
$cell.Value = "'158.27"
$cell.Style = 'Normal'
$ws.Range('E26').Value = '  +0.96%  '
$cell = $ws.Range('D27')
$cell.Value = "'19.04"
$cell.Style = 'Normal'
$ws.Range('E27').Value = '  -0.66%  '
$cell = $ws.Range('D28')
$cell.Value = "'5.636"
$cell.Style = 'Normal'
$ws.Range('E28').Value = '  +0.81%  '
$cell = $ws.Range('D29')
$cell.Value = "'117.66"
$cell.Style = 'Normal'
$ws.Range('E29').Value = '  -0.11%  '
$cell = $ws.Range('D30')
$cell.Value = "'1.842"
$cell.Style = 'Normal'
$cell = $ws.Range('D31')
$cell.Value = "'0.09298"
$cell.Style = 'Normal'
$ws.Range('E31').Value = '  +0.32%  '
$ws.Range('E32').Value = '  +0.46%  '
$cell = $ws.Range('D33')
$cell.Value = "'5.116"
$cell.Style = 'Normal'
$ws.Range('E33').Value = '  +0.27%  '
$cell = $ws.Range('D34')
$cell.Value = "'1.249"
$cell.Style = 'Normal'
$ws.Range('E34').Value = '  -0.07%  '
$cell = $ws.Range('D35')
$cell.Value = "'3.017"
$cell.Style = 'Normal'
$ws.Range('E35').Value = '  +0.09%  '
$cell = $ws.Range('D36')
$cell.Value = "'0.05710"
$cell.Style = 'Normal'
$ws.Range('E36').Value = '  +0.01%  '
$cell = $ws.Range('D37')
$cell.Value = "'1.153"
$cell.Style = 'Normal'
$ws.Range('E37').Value = '  +0.76%  '
$ws.Range('E38').Value = '  +0.10%  '
$cell = $ws.Range('D39')
$cell.Value = "'0.02056"
$cell.Style = 'Normal'
$ws.Range('E39').Value = '  +0.99%  '
$cell = $ws.Range('D40')
$cell.Value = "'3.102"
$cell.Style = 'Normal'
$ws.Range('E40').Value = '  +13.55%  '
$cell = $ws.Range('D41')
$cell.Value = "'7.477"
$cell.Style = 'Normal'
$ws.Range('E41').Value = '  -0.19%  '
$cell = $ws.Range('D42')
$cell.Value = "'0.5521"
$cell.Style = 'Normal'
$ws.Range('E42').Value = '  +0.11%  '
$cell = $ws.Range('D43')
$cell.Value = "'0.1757"
$cell.Style = 'Normal'
$ws.Range('E43').Value = '  +0.03%  '
$cell = $ws.Range('D44')
$cell.Value = "'9.375"
$cell.Style = 'Normal'
$ws.Range('E44').Value = '  +0.49%  '
$cell = $ws.Range('D45')
$cell.Value = "'0.000002844"
$cell.Style = 'Normal'
$ws.Range('E45').Value = '  +8.75%  '
$cell = $ws.Range('D46')
$cell.Value = "'2.186"
$cell.Style = 'Normal'
$ws.Range('E46').Value = '  +4.68%  '
$cell = $ws.Range('D47')
$cell.Value = "'0.5188"
$cell.Style = 'Normal'
$ws.Range('E47').Value = '  -0.29%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$cell = $ws.Range('D48')
$cell.Value = "'11.25"
$cell.Style = 'Normal'
$ws.Range('E48').Value = '  -0.38%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$cell = $ws.Range('D49')
$cell.Value = "'0.06937"
$cell.Style = 'Normal'
$ws.Range('E49').Value = '  +1.80%  '
$cell = $ws.Range('D50')
$cell.Value = "'111.11"
$cell.Style = 'Normal'
$ws.Range('E50').Value = '  -0.36%  '
$cell = $ws.Range('D51')
$cell.Value = "'1.766"
$cell.Style = 'Normal'
$ws.Range('E51').Value = '  -0.46%  '
